{"js": "const body = context.document.body;\n\n// Ordered list of (old, new) text replacements taken from the diff.\n// Every old value is unique and occurs exactly once in the document,\n// and no new value collides with any other old value, so a plain\n// search-and-replace per pair is unambiguous and order independent.\nconst pairs = [\n  [\"34\u00d726=\", \"81\u00d782=\"],\n  [\"69\u00d774=\", \"54\u00d737=\"],\n  [\"80\u00d735=\", \"29\u00d736=\"],\n  [\"27\u00d755=\", \"14\u00d727=\"],\n  [\"47\u00d781=\", \"13\u00d754=\"],\n  [\"42\u00d795=\", \"48\u00d762=\"],\n  [\"89\u00d741=\", \"39\u00d771=\"],\n  [\"98\u00d757=\", \"99\u00d789=\"],\n  [\"70\u00d789=\", \"44\u00d712=\"],\n  [\"53\u00d7100=\", \"14\u00d740=\"],\n  [\"41\u00d742=\", \"22\u00d735=\"],\n  [\"51\u00d770=\", \"14\u00d721=\"],\n  [\"70\u00d715=\", \"64\u00d748=\"],\n  [\"91\u00d725=\", \"85\u00d799=\"],\n  [\"35\u00d733=\", \"82\u00d714=\"],\n  [\"75\u00d784=\", \"99\u00d764=\"],\n  [\"32\u00d759=\", \"17\u00d740=\"],\n  [\"25\u00d755=\", \"41\u00d798=\"],\n  [\"99\u00d769=\", \"63\u00d748=\"],\n  [\"80\u00d776=\", \"98\u00d727=\"],\n  [\"84\u00d772=\", \"76\u00d752=\"],\n  [\"52\u00d735=\", \"46\u00d744=\"],\n  [\"66\u00d724=\", \"53\u00d732=\"],\n  [\"44\u00d725=\", \"89\u00d718=\"],\n  [\"62\u00d770=\", \"74\u00d749=\"],\n  [\"10\u00d741=\", \"28\u00d784=\"],\n  [\"37\u00d765=\", \"29\u00d794=\"],\n  [\"48\u00d711=\", \"25\u00d716=\"],\n  [\"84\u00d725=\", \"27\u00d798=\"],\n  [\"75\u00d778=\", \"69\u00d781=\"],\n  [\"86\u00d792=\", \"48\u00d784=\"],\n  [\"57\u00d781=\", \"42\u00d769=\"],\n  [\"93\u00d734=\", \"74\u00d738=\"],\n  [\"27\u00d742=\", \"60\u00d778=\"],\n  [\"97\u00d760=\", \"21\u00d727=\"],\n  [\"43\u00d744=\", \"82\u00d767=\"],\n  [\"80\u00d710=\", \"25\u00d758=\"],\n  [\"47\u00d712=\", \"49\u00d799=\"],\n  [\"59\u00d746=\", \"33\u00d786=\"],\n  [\"69\u00d754=\", \"72\u00d798=\"],\n  [\"91\u00d723=\", \"89\u00d766=\"],\n  [\"68\u00d734=\", \"25\u00d733=\"],\n  [\"90\u00d746=\", \"86\u00d761=\"],\n  [\"24\u00d789=\", \"22\u00d738=\"],\n  [\"51\u00d792=\", \"98\u00d797=\"],\n  [\"70\u00d7100=\", \"44\u00d782=\"],\n  [\"21\u00d789=\", \"78\u00d744=\"],\n  [\"83\u00d770=\", \"32\u00d768=\"],\n  [\"91\u00d775=\", \"97\u00d774=\"],\n  [\"36\u00d797=\", \"28\u00d757=\"],\n  [\"86\u00d721=\", \"43\u00d719=\"],\n  [\"98\u00d7100=\", \"79\u00d731=\"],\n  [\"86\u00d784=\", \"15\u00d716=\"],\n  [\"97\u00d755=\", \"26\u00d729=\"],\n  [\"82\u00d788=\", \"43\u00d780=\"],\n  [\"33\u00d760=\", \"81\u00d774=\"],\n  [\"94\u00d721=\", \"81\u00d786=\"],\n  [\"60\u00d789=\", \"80\u00d758=\"],\n  [\"74\u00d739=\", \"54\u00d785=\"],\n  [\"34\u00d740=\", \"29\u00d768=\"],\n  [\"39\u00d764=\", \"85\u00d757=\"],\n  [\"21\u00d720=\", \"20\u00d764=\"],\n  [\"72\u00d740=\", \"92\u00d795=\"],\n  [\"45\u00d729=\", \"77\u00d752=\"],\n  [\"34\u00d763=\", \"59\u00d795=\"],\n  [\"30\u00d796=\", \"18\u00d735=\"],\n  [\"45\u00d785=\", \"20\u00d769=\"],\n  [\"33\u00d736=\", \"51\u00d781=\"],\n  [\"53\u00d792=\", \"18\u00d754=\"],\n  [\"74\u00d735=\", \"65\u00d749=\"],\n  [\"82\u00d728=\", \"55\u00d790=\"],\n  [\"27\u00d790=\", \"33\u00d731=\"],\n  [\"76\u00d715=\", \"74\u00d768=\"],\n  [\"41\u00d799=\", \"16\u00d741=\"],\n  [\"80\u00d766=\", \"24\u00d721=\"],\n  [\"78\u00d761=\", \"15\u00d778=\"],\n  [\"16\u00d746=\", \"62\u00d729=\"],\n  [\"63\u00d766=\", \"47\u00d756=\"],\n  [\"81\u00d760=\", \"87\u00d776=\"],\n  [\"52\u00d721=\", \"88\u00d769=\"],\n  [\"65\u00d712=\", \"12\u00d720=\"],\n  [\"35\u00d753=\", \"15\u00d795=\"],\n  [\"13\u00d782=\", \"24\u00d717=\"],\n  [\"71\u00d784=\", \"44\u00d758=\"],\n  [\"35\u00d768=\", \"91\u00d737=\"],\n  [\"23\u00d762=\", \"17\u00d756=\"],\n  [\"49\u00d743=\", \"26\u00d784=\"],\n  [\"29\u00d770=\", \"84\u00d730=\"],\n  [\"80\u00d775=\", \"29\u00d791=\"],\n  [\"78\u00d754=\", \"81\u00d735=\"],\n  [\"57\u00d770=\", \"54\u00d714=\"],\n  [\"65\u00d769=\", \"38\u00d738=\"],\n  [\"61\u00d770=\", \"16\u00d733=\"],\n  [\"13\u00d748=\", \"60\u00d757=\"],\n  [\"12\u00d770=\", \"12\u00d719=\"],\n  [\"69\u00d785=\", \"43\u00d755=\"],\n  [\"100\u00d711=\", \"66\u00d740=\"],\n  [\"95\u00d794=\", \"55\u00d799=\"],\n  [\"73\u00d717=\", \"26\u00d770=\"],\n  [\"55\u00d757=\", \"24\u00d730=\"],\n];\n\n// Phase 1: issue a search for every old value and queue the loads.\nconst searches = pairs.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nsearches.forEach((r) => r.load(\"items\"));\nawait context.sync();\n\n// Phase 2: replace the text of each found range in place. Using\n// Range.insertText(text, \"Replace\") on the *found range* (rather than\n// rewriting the whole paragraph/cell) preserves the existing run\n// formatting (font, size, etc.) and paragraph properties untouched.\nsearches.forEach((results, i) => {\n  const [oldText, newText] = pairs[i];\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n});\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Ordered list of (old, new) text replacements taken from the diff.\n# Every old value is unique and occurs exactly once in the document,\n# and no new value collides with any other old value, so a plain\n# Find/Replace per pair is unambiguous and order independent.\n$pairs = @(\n  @('34\u00d726=', '81\u00d782='),\n  @('69\u00d774=', '54\u00d737='),\n  @('80\u00d735=', '29\u00d736='),\n  @('27\u00d755=', '14\u00d727='),\n  @('47\u00d781=', '13\u00d754='),\n  @('42\u00d795=', '48\u00d762='),\n  @('89\u00d741=', '39\u00d771='),\n  @('98\u00d757=', '99\u00d789='),\n  @('70\u00d789=', '44\u00d712='),\n  @('53\u00d7100=', '14\u00d740='),\n  @('41\u00d742=', '22\u00d735='),\n  @('51\u00d770=', '14\u00d721='),\n  @('70\u00d715=', '64\u00d748='),\n  @('91\u00d725=', '85\u00d799='),\n  @('35\u00d733=', '82\u00d714='),\n  @('75\u00d784=', '99\u00d764='),\n  @('32\u00d759=', '17\u00d740='),\n  @('25\u00d755=', '41\u00d798='),\n  @('99\u00d769=', '63\u00d748='),\n  @('80\u00d776=', '98\u00d727='),\n  @('84\u00d772=', '76\u00d752='),\n  @('52\u00d735=', '46\u00d744='),\n  @('66\u00d724=', '53\u00d732='),\n  @('44\u00d725=', '89\u00d718='),\n  @('62\u00d770=', '74\u00d749='),\n  @('10\u00d741=', '28\u00d784='),\n  @('37\u00d765=', '29\u00d794='),\n  @('48\u00d711=', '25\u00d716='),\n  @('84\u00d725=', '27\u00d798='),\n  @('75\u00d778=', '69\u00d781='),\n  @('86\u00d792=', '48\u00d784='),\n  @('57\u00d781=', '42\u00d769='),\n  @('93\u00d734=', '74\u00d738='),\n  @('27\u00d742=', '60\u00d778='),\n  @('97\u00d760=', '21\u00d727='),\n  @('43\u00d744=', '82\u00d767='),\n  @('80\u00d710=', '25\u00d758='),\n  @('47\u00d712=', '49\u00d799='),\n  @('59\u00d746=', '33\u00d786='),\n  @('69\u00d754=', '72\u00d798='),\n  @('91\u00d723=', '89\u00d766='),\n  @('68\u00d734=', '25\u00d733='),\n  @('90\u00d746=', '86\u00d761='),\n  @('24\u00d789=', '22\u00d738='),\n  @('51\u00d792=', '98\u00d797='),\n  @('70\u00d7100=', '44\u00d782='),\n  @('21\u00d789=', '78\u00d744='),\n  @('83\u00d770=', '32\u00d768='),\n  @('91\u00d775=', '97\u00d774='),\n  @('36\u00d797=', '28\u00d757='),\n  @('86\u00d721=', '43\u00d719='),\n  @('98\u00d7100=', '79\u00d731='),\n  @('86\u00d784=', '15\u00d716='),\n  @('97\u00d755=', '26\u00d729='),\n  @('82\u00d788=', '43\u00d780='),\n  @('33\u00d760=', '81\u00d774='),\n  @('94\u00d721=', '81\u00d786='),\n  @('60\u00d789=', '80\u00d758='),\n  @('74\u00d739=', '54\u00d785='),\n  @('34\u00d740=', '29\u00d768='),\n  @('39\u00d764=', '85\u00d757='),\n  @('21\u00d720=', '20\u00d764='),\n  @('72\u00d740=', '92\u00d795='),\n  @('45\u00d729=', '77\u00d752='),\n  @('34\u00d763=', '59\u00d795='),\n  @('30\u00d796=', '18\u00d735='),\n  @('45\u00d785=', '20\u00d769='),\n  @('33\u00d736=', '51\u00d781='),\n  @('53\u00d792=', '18\u00d754='),\n  @('74\u00d735=', '65\u00d749='),\n  @('82\u00d728=', '55\u00d790='),\n  @('27\u00d790=', '33\u00d731='),\n  @('76\u00d715=', '74\u00d768='),\n  @('41\u00d799=', '16\u00d741='),\n  @('80\u00d766=', '24\u00d721='),\n  @('78\u00d761=', '15\u00d778='),\n  @('16\u00d746=', '62\u00d729='),\n  @('63\u00d766=', '47\u00d756='),\n  @('81\u00d760=', '87\u00d776='),\n  @('52\u00d721=', '88\u00d769='),\n  @('65\u00d712=', '12\u00d720='),\n  @('35\u00d753=', '15\u00d795='),\n  @('13\u00d782=', '24\u00d717='),\n  @('71\u00d784=', '44\u00d758='),\n  @('35\u00d768=', '91\u00d737='),\n  @('23\u00d762=', '17\u00d756='),\n  @('49\u00d743=', '26\u00d784='),\n  @('29\u00d770=', '84\u00d730='),\n  @('80\u00d775=', '29\u00d791='),\n  @('78\u00d754=', '81\u00d735='),\n  @('57\u00d770=', '54\u00d714='),\n  @('65\u00d769=', '38\u00d738='),\n  @('61\u00d770=', '16\u00d733='),\n  @('13\u00d748=', '60\u00d757='),\n  @('12\u00d770=', '12\u00d719='),\n  @('69\u00d785=', '43\u00d755='),\n  @('100\u00d711=', '66\u00d740='),\n  @('95\u00d794=', '55\u00d799='),\n  @('73\u00d717=', '26\u00d770='),\n  @('55\u00d757=', '24\u00d730='),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  # wdFindContinue (1) replaces just the one match per Execute call;\n  # Forward=True, Wrap=wdFindContinue keeps the search scoped to $d.Content.\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n\n  $found = $find.Execute(\n    [ref]$oldText,    # FindText\n    [ref]$false,      # MatchCase\n    [ref]$false,      # MatchWholeWord\n    [ref]$false,      # MatchWildcards\n    [ref]$false,      # MatchSoundsLike\n    [ref]$false,      # MatchAllWordForms\n    [ref]$true,       # Forward\n    1,                # Wrap (wdFindContinue)\n    [ref]$false,      # Format\n    $newText,         # ReplaceWith\n    2                 # Replace (wdReplaceOne) - replace this single match only\n  )\n\n  if (-not $found) {\n    throw \"Not found: $oldText\"\n  }\n}\n"}
